# Auto-generated Excel COM-interop script
# Applies the numeric cell updates described by the commit diff
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1357.25
$ws.Range("I12").Value = 793.3333
$ws.Range("J12").Value = 1695.6
$ws.Range("K12").Value = 793.3333
$ws.Range("L12").Value = 1695.6
$ws.Range("M12").Value = -623.3333
$ws.Range("N12").Value = -2035.6
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = ""
$ws.Range("H69").Value = 17499
$ws.Range("I69").Value = 17499
$ws.Range("K69").Value = 52497
$ws.Range("M69").Value = -51623
$ws.Range("H72").Value = 17499
$ws.Range("I72").Value = 17499
$ws.Range("K72").Value = 157491
$ws.Range("M72").Value = -153123
$ws.Range("H80").Value = 167513.25
$ws.Range("I80").Value = 334254.84
$ws.Range("J80").Value = 771.6667
$ws.Range("K80").Value = 1002764.52
$ws.Range("L80").Value = 2315.0001
$ws.Range("M80").Value = -1001766.52
$ws.Range("N80").Value = -4311.0001
$ws.Range("H83").Value = 167513.25
$ws.Range("I83").Value = 334254.84
$ws.Range("J83").Value = 771.6667
$ws.Range("K83").Value = 3008293.56
$ws.Range("L83").Value = 6945.0003
$ws.Range("M83").Value = -3003301.56
$ws.Range("N83").Value = -16929.0003
$ws.Range("H86").Value = 30027.375
$ws.Range("I86").Value = 3687.6
$ws.Range("K86").Value = 3687.6
$ws.Range("M86").Value = -2564.6
$ws.Range("H89").Value = 30027.375
$ws.Range("I89").Value = 3687.6
$ws.Range("K89").Value = 18438
$ws.Range("M89").Value = -12822
$ws.Range("H113").Value = 7780.75
$ws.Range("I113").Value = 7696.9
$ws.Range("J113").Value = 8200
$ws.Range("K113").Value = 7696.9
$ws.Range("L113").Value = 8200
$ws.Range("M113").Value = -4442.9
$ws.Range("N113").Value = -14708

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 7514.9736
$ws.Range("I137").Value = 8058.353
$ws.Range("K137").Value = 24175.059
$ws.Range("M137").Value = -21625.059
$ws.Range("H138").Value = 3673.451
$ws.Range("I138").Value = 1942.7273
$ws.Range("J138").Value = 4149.4
$ws.Range("K138").Value = 5828.1819
$ws.Range("L138").Value = 12448.2
$ws.Range("M138").Value = -688.1818999999996
$ws.Range("N138").Value = -22728.2
$ws.Range("H2").Value = 75762
$ws.Range("I2").Value = 12318
$ws.Range("K2").Value = 12318
$ws.Range("M2").Value = -12205
$ws.Range("H44").Value = 97999.336
$ws.Range("J44").Value = 97999.336
$ws.Range("L44").Value = 97999.336
$ws.Range("N44").Value = -98975.336
$ws.Range("H61").Value = 3968.261
$ws.Range("I61").Value = 3213.543
$ws.Range("K61").Value = 3213.543
$ws.Range("M61").Value = -3001.543
$ws.Range("H74").Value = 2739.6206
$ws.Range("I74").Value = 1570.7916
$ws.Range("K74").Value = 1570.7916
$ws.Range("M74").Value = -696.7916
$ws.Range("H77").Value = 2739.6206
$ws.Range("I77").Value = 1570.7916
$ws.Range("K77").Value = 7853.958000000001
$ws.Range("M77").Value = -3485.958000000001
$ws.Range("H110").Value = 3433.1333
$ws.Range("J110").Value = 4800
$ws.Range("L110").Value = 4800
$ws.Range("N110").Value = -8890
$ws.Range("H116").Value = 75762
$ws.Range("I116").Value = 12318
$ws.Range("K116").Value = 12318
$ws.Range("M116").Value = -10024

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 3968.261
$ws.Range("I136").Value = 3213.543
$ws.Range("K136").Value = 9640.629000000001
$ws.Range("M136").Value = -7090.629000000001
$ws.Range("H3").Value = 75762
$ws.Range("I3").Value = 12318
$ws.Range("K3").Value = 12318
$ws.Range("M3").Value = -12204
$ws.Range("H86").Value = 6026.387
$ws.Range("I86").Value = 4983.773
$ws.Range("J86").Value = 8575
$ws.Range("K86").Value = 4983.773
$ws.Range("L86").Value = 8575
$ws.Range("M86").Value = -3860.773
$ws.Range("N86").Value = -10821
$ws.Range("H89").Value = 6026.387
$ws.Range("I89").Value = 4983.773
$ws.Range("J89").Value = 8575
$ws.Range("K89").Value = 24918.865
$ws.Range("L89").Value = 42875
$ws.Range("M89").Value = -19302.865
$ws.Range("N89").Value = -54107
$ws.Range("H107").Value = 3812.7083
$ws.Range("J107").Value = 4671.4287
$ws.Range("L107").Value = 4671.4287
$ws.Range("N107").Value = -8511.4287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 10002
$ws.Range("I3").Value = 10002
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 10002
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -9889
$ws.Range("N3").Value = ""
$ws.Range("H105").Value = 9849.333000000001
$ws.Range("I105").Value = 13524
$ws.Range("K105").Value = 13524
$ws.Range("M105").Value = -11777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 24536450
$ws.Range("I4").Value = 3890114
$ws.Range("K4").Value = 11670342
$ws.Range("M4").Value = -11670230

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 10000
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = ""
$ws.Range("H114").Value = 60000
$ws.Range("J114").Value = 60000
$ws.Range("L114").Value = 60000
$ws.Range("N114").Value = -68678

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4294.4614
$ws.Range("I132").Value = 4980.8887
$ws.Range("J132").Value = 2750
$ws.Range("K132").Value = 14942.6661
$ws.Range("L132").Value = 8250
$ws.Range("M132").Value = -12412.6661
$ws.Range("N132").Value = -13310
$ws.Range("H22").Value = 29486
$ws.Range("I22").Value = 29486
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 29486
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -29191
$ws.Range("N22").Value = ""
$ws.Range("H27").Value = 29486
$ws.Range("I27").Value = 29486
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 29486
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -29379
$ws.Range("N27").Value = ""
$ws.Range("H38").Value = 38570.43
$ws.Range("I38").Value = 20030
$ws.Range("J38").Value = 41660.5
$ws.Range("K38").Value = 20030
$ws.Range("L38").Value = 41660.5
$ws.Range("M38").Value = -19620
$ws.Range("N38").Value = -42480.5
$ws.Range("H40").Value = 56408.445
$ws.Range("I40").Value = 63955.145
$ws.Range("K40").Value = 63955.145
$ws.Range("M40").Value = -63819.145
$ws.Range("H46").Value = 2819.68
$ws.Range("I46").Value = 1736.1818
$ws.Range("K46").Value = 1736.1818
$ws.Range("M46").Value = -1548.1818
$ws.Range("H55").Value = 1478.9286
$ws.Range("I55").Value = 204.25
$ws.Range("K55").Value = 204.25
$ws.Range("M55").Value = -31.25
$ws.Range("H68").Value = 5468.222
$ws.Range("J68").Value = 6729.6924
$ws.Range("L68").Value = 6729.6924
$ws.Range("N68").Value = -8227.6924
$ws.Range("H71").Value = 5468.222
$ws.Range("J71").Value = 6729.6924
$ws.Range("L71").Value = 33648.462
$ws.Range("N71").Value = -41136.462
$ws.Range("H93").Value = 14776.23
$ws.Range("I93").Value = 14177.111
$ws.Range("J93").Value = 16124.25
$ws.Range("K93").Value = 14177.111
$ws.Range("L93").Value = 16124.25
$ws.Range("M93").Value = -12929.111
$ws.Range("N93").Value = -18620.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3557.8
$ws.Range("I136").Value = 1282.1818
$ws.Range("J136").Value = 7975.1763
$ws.Range("K136").Value = 3846.5454
$ws.Range("L136").Value = 23925.5289
$ws.Range("M136").Value = -1296.5454
$ws.Range("N136").Value = -29025.5289
$ws.Range("H11").Value = 2000
$ws.Range("I11").Value = 2000
$ws.Range("K11").Value = 2000
$ws.Range("M11").Value = -1858
$ws.Range("H62").Value = 148806.27
$ws.Range("J62").Value = 3784.9285
$ws.Range("L62").Value = 3784.9285
$ws.Range("N62").Value = -5032.9285
$ws.Range("H65").Value = 148806.27
$ws.Range("J65").Value = 3784.9285
$ws.Range("L65").Value = 18924.6425
$ws.Range("N65").Value = -25164.6425
$ws.Range("H81").Value = 16163.818
$ws.Range("I81").Value = 16980.2
$ws.Range("J81").Value = 8000
$ws.Range("K81").Value = 33960.4
$ws.Range("L81").Value = 16000
$ws.Range("M81").Value = -32899.4
$ws.Range("N81").Value = -18122
$ws.Range("H84").Value = 16163.818
$ws.Range("I84").Value = 16980.2
$ws.Range("J84").Value = 8000
$ws.Range("K84").Value = 169802
$ws.Range("L84").Value = 80000
$ws.Range("M84").Value = -164498
$ws.Range("N84").Value = -90608
$ws.Range("H126").Value = 18447.629
$ws.Range("J126").Value = 6456.2
$ws.Range("L126").Value = 19368.6
$ws.Range("N126").Value = -24308.6
$ws.Range("H132").Value = 7458.4326
$ws.Range("I132").Value = 7707.3965
$ws.Range("J132").Value = 6555.9375
$ws.Range("K132").Value = 23122.1895
$ws.Range("L132").Value = 19667.8125
$ws.Range("M132").Value = -20592.1895
$ws.Range("N132").Value = -24727.8125
$ws.Range("H136").Value = 1189.44
$ws.Range("I136").Value = 926.5
$ws.Range("K136").Value = 2779.5
$ws.Range("M136").Value = -229.5
